# Rename the worksheet from "example_24_well" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Update the active selection on the sheet to G27 (single cell)
$ws.Range("G27").Select()
